# Madagascar (Madagascar/MDG) localisation of machine_spec.xlsx:
# - re-language the single data row from English ("eng") to French ("fra")
# - refresh header/body styling (borders, fonts, wrap text, alignment)
#   to match the look produced when the workbook is re-saved by Excel
# - resize columns, set explicit row heights, tidy page margins

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Data: translate the single data row to French (Madagascar locale)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "fra"
$ws.Range("C2").Value = "Machine virtuelle résidente"
$ws.Range("H2").Value = "Machine virtuelle résidente"

# ---------------------------------------------------------------------
# 2. Column widths
# ---------------------------------------------------------------------
$offset = 5.0/6.0
$ws.Columns.Item(1).ColumnWidth = 19.453125 - $offset
$ws.Columns.Item(2).ColumnWidth = 18 - $offset
$ws.Columns.Item(3).ColumnWidth = 36.54296875 - $offset
$ws.Columns.Item(4).ColumnWidth = 38.6328125 - $offset
$ws.Columns.Item(5).ColumnWidth = 18.7265625 - $offset
$ws.Columns.Item(6).ColumnWidth = 12.90625 - $offset
$ws.Columns.Item(7).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(8).ColumnWidth = 41.54296875 - $offset
$ws.Columns.Item(9).ColumnWidth = $ws.StandardWidth

# ---------------------------------------------------------------------
# 3. Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 28.5
$ws.Rows.Item(2).RowHeight = 44

# ---------------------------------------------------------------------
# 4. Header row (row 1) formatting
# ---------------------------------------------------------------------
$headerAll = $ws.Range("A1:I1")
$headerAll.Font.Name = "Cambria"
$headerAll.Font.Bold = $true
$headerAll.Font.ThemeColor = 1
$headerAll.HorizontalAlignment = -4108   # xlCenter
$headerAll.VerticalAlignment = -4160     # xlTop
$headerAll.WrapText = $true
$headerAll.NumberFormat = "General"

# A1 gets a uniform medium black border on all sides
$a1 = $ws.Range("A1")
$a1.Borders.Weight = -4138               # xlMedium
$a1.Borders.Color = 0

# B1:I1 get a medium black border on top/right/bottom and a light-grey
# medium border on the left (matches the vertical "light separator" look)
$restHeader = $ws.Range("B1:I1")
$restHeader.Borders.Weight = -4138
foreach ($idx in 8,9,10) {
    $restHeader.Borders.Item($idx).Color = 0
}
$restHeader.Borders.Item(7).Color = 13421772   # CCCCCC

# ---------------------------------------------------------------------
# 5. Data row (row 2) formatting
# ---------------------------------------------------------------------
$dataAll = $ws.Range("A2:I2")
$dataAll.Font.Name = "Calibri"
$dataAll.Font.Bold = $false
$dataAll.Font.ThemeColor = 1
$dataAll.WrapText = $true
$dataAll.Borders.Weight = -4138          # xlMedium, light grey all round
$dataAll.Borders.Color = 13421772

# G2 (is boolean? no: min_driver_ver numeric placeholder) -> right aligned
$ws.Range("G2").HorizontalAlignment = -4152   # xlRight

# I2 (is_active boolean) -> center aligned, plain General number format
$i2 = $ws.Range("I2")
$i2.HorizontalAlignment = -4108               # xlCenter
$i2.NumberFormat = "General"

# ---------------------------------------------------------------------
# 6. Misc workbook / sheet view bits
# ---------------------------------------------------------------------
[void]$ws.Range("D10").Select()

$ps = $ws.PageSetup
$ps.LeftMargin = 50.4
$ps.RightMargin = 50.4
$ps.TopMargin = 54
$ps.BottomMargin = 54
$ps.HeaderMargin = 21.6
$ps.FooterMargin = 21.6

$wb.ChartDataPointTrack = $true
